$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")
$ws.Activate()

# --- Update the reporting-period / validation dates for rows 8-10 ---
# Columns: B = start date, C = end date, J = validation date, K = update date
# Q3 2021 (Jul-Sep) values move forward to Q4 2021 (Oct-Dec) values.
$rows = 8,9,10
foreach ($r in $rows) {
    $ws.Cells.Item($r, 2).Value = 44470   # B -> 2021-10-01
    $ws.Cells.Item($r, 3).Value = 44561   # C -> 2021-12-31
    $ws.Cells.Item($r, 10).Value = 44571  # J -> 2022-01-10
    $ws.Cells.Item($r, 11).Value = 44571  # K -> 2022-01-10
}

# --- Narrow column E slightly (41.57 -> ~40.14 characters) ---
$ws.Columns.Item(5).ColumnWidth = 39.3

# --- Update the view: scroll so row/col B2 is the top-left visible cell, ---
# --- and move the active selection to H11 ---
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("H11").Select()
